$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.652.46'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').Value = '2.127.23'
$ws.Range('E3').Value = '  +1.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.013'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '352.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.011'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.74%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5280'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4550'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09105'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.91%  '
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.66'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('D13').Value = '2.132.26'
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.875'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.131'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '102.43'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.03%  '
$ws.Range('E17').Value = '  +3.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.013'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06720'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.50'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.011'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.363'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').Value = '30.738.23'
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('E24').Value = '  +3.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.386'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('D26').Value = '2.384.35'
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.571'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '164.83'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '136.48'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.202'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1086'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.671'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.395'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.022'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.175'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.32'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02655'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06892'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2325'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.58'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6926'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.279'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.76'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.33%  '
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6465'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.33%  '
$ws.Range('E47').Value = '  +2.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000365'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.59%  '
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3427'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '83.16'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.40%  '
